$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell in the Price (D) / Volume(1h) (E) columns holds text (inline string)
# in the workbook. Forcing the cell to Text format before assigning, then
# restoring the style to Normal afterwards, guarantees the value is written
# back as a string (preserving things like leading/trailing zeros, e.g.
# "1.00", "4.20") instead of being auto-coerced into a number by Excel,
# while leaving the cell formatting identical to before the edit.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.055.82"
Set-TextValue $ws.Range("E2") "  +1.44%  "
Set-TextValue $ws.Range("D3") "1.643.97"
Set-TextValue $ws.Range("E4") "  +0.18%  "
Set-TextValue $ws.Range("D5") "215.63"
Set-TextValue $ws.Range("D6") "0.521"
Set-TextValue $ws.Range("E6") "  +1.10%  "
Set-TextValue $ws.Range("E7") "  +0.16%  "
Set-TextValue $ws.Range("D8") "29.41"
Set-TextValue $ws.Range("E8") "  +5.37%  "
Set-TextValue $ws.Range("E9") "  +3.53%  "
Set-TextValue $ws.Range("D10") "0.0614"
Set-TextValue $ws.Range("E10") "  +1.95%  "
Set-TextValue $ws.Range("D11") "0.0917"
Set-TextValue $ws.Range("E11") "  +0.56%  "
Set-TextValue $ws.Range("D12") "1.879.62"
Set-TextValue $ws.Range("E12") "  +2.65%  "
Set-TextValue $ws.Range("D13") "1.646.13"
Set-TextValue $ws.Range("E13") "  +2.77%  "
Set-TextValue $ws.Range("D14") "0.573"
Set-TextValue $ws.Range("E14") "  +5.35%  "
Set-TextValue $ws.Range("D15") "9.48"
Set-TextValue $ws.Range("E15") "  +22.18%  "
Set-TextValue $ws.Range("E16") "  +4.74%  "
Set-TextValue $ws.Range("D17") "30.100.54"
Set-TextValue $ws.Range("E17") "  +1.56%  "
Set-TextValue $ws.Range("E18") "  +1.74%  "
Set-TextValue $ws.Range("D19") "247.76"
Set-TextValue $ws.Range("E19") "  +2.21%  "
Set-TextValue $ws.Range("D20") "0.0₃0710"
Set-TextValue $ws.Range("E20") "  +1.99%  "
Set-TextValue $ws.Range("D21") "1.00"
Set-TextValue $ws.Range("E21") "  +0.12%  "
Set-TextValue $ws.Range("D22") "10.02"
Set-TextValue $ws.Range("E22") "  +6.80%  "
Set-TextValue $ws.Range("D23") "4.20"
Set-TextValue $ws.Range("E23") "  +4.70%  "
Set-TextValue $ws.Range("D24") "2.16"
Set-TextValue $ws.Range("E24") "  +2.59%  "
Set-TextValue $ws.Range("D25") "158.91"
Set-TextValue $ws.Range("E25") "  +2.48%  "
Set-TextValue $ws.Range("D26") "15.76"
Set-TextValue $ws.Range("E26") "  +2.17%  "
Set-TextValue $ws.Range("E27") "  +2.57%  "
Set-TextValue $ws.Range("D28") "6.67"
Set-TextValue $ws.Range("E28") "  +3.87%  "
Set-TextValue $ws.Range("E29") "  +0.09%  "
Set-TextValue $ws.Range("E30") "  +2.63%  "
Set-TextValue $ws.Range("E31") "  +6.13%  "
Set-TextValue $ws.Range("E32") "  +6.32%  "
Set-TextValue $ws.Range("D33") "3.21"
Set-TextValue $ws.Range("E33") "  +0.82%  "
Set-TextValue $ws.Range("D34") "1.441.54"
Set-TextValue $ws.Range("E34") "  +1.14%  "
Set-TextValue $ws.Range("E35") "  +7.63%  "
Set-TextValue $ws.Range("E36") "  +1.87%  "
Set-TextValue $ws.Range("E37") "  -0.92%  "
Set-TextValue $ws.Range("D38") "78.18"
Set-TextValue $ws.Range("E38") "  +18.16%  "
Set-TextValue $ws.Range("E39") "  +2.01%  "
Set-TextValue $ws.Range("E40") "  -0.10%  "
Set-TextValue $ws.Range("E41") "  +2.61%  "
Set-TextValue $ws.Range("E42") "  +2.72%  "
Set-TextValue $ws.Range("D43") "0.846"
Set-TextValue $ws.Range("E43") "  +3.57%  "
Set-TextValue $ws.Range("D44") "55.49"
Set-TextValue $ws.Range("E44") "  -3.18%  "
Set-TextValue $ws.Range("D45") "0.0499"
Set-TextValue $ws.Range("E45") "  +0.78%  "
Set-TextValue $ws.Range("E46") "  +6.25%  "
Set-TextValue $ws.Range("E47") "  +0.11%  "
Set-TextValue $ws.Range("D48") "5.39"
Set-TextValue $ws.Range("E48") "  +0.86%  "
Set-TextValue $ws.Range("D49") "1.785.95"
Set-TextValue $ws.Range("E49") "  +2.67%  "
Set-TextValue $ws.Range("E50") "  +11.59%  "
Set-TextValue $ws.Range("D51") "90.41"
Set-TextValue $ws.Range("E51") "  +4.10%  "
